$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.710.67'
$ws.Range('E2').Value = '  +0.33%  '
$ws.Range('D3').Value = '1.600.58'
$ws.Range('E3').Value = '  +0.25%  '
$ws.Range('E4').Value = '  +0.31%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.54'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.03%  '
$ws.Range('E6').Value = '  -0.56%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.01'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.27%  '
$ws.Range('E9').Value = '  +0.63%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.56'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.35%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0843'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.71%  '
$ws.Range('D12').Value = '1.825.60'
$ws.Range('E12').Value = '  +0.28%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.05'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.49%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.563.02'
$ws.Range('E14').Value = '  -1.11%  '
$ws.Range('E15').Value = '  +0.24%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.32'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.06%  '
$ws.Range('D17').Value = '26.683.34'
$ws.Range('D18').Value = '0.0₃0757'
$ws.Range('E18').Value = '  +3.40%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.25'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.18%  '
$ws.Range('E20').Value = '  +0.31%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '209.40'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.29%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.29'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.47%  '
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('E24').Value = '  +0.55%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '142.84'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.80%  '
$ws.Range('E26').Value = '  +0.20%  '
$ws.Range('E27').Value = '  -0.14%  '
$ws.Range('E28').Value = '  -0.15%  '
$ws.Range('E29').Value = '  +0.57%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.16'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.07%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.26'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.77%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.98'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.35%  '
$ws.Range('D34').Value = '1.294.39'
$ws.Range('E34').Value = '  +0.95%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.621'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.27%  '
$ws.Range('E36').Value = '  +1.07%  '
$ws.Range('E37').Value = '  +0.28%  '
$ws.Range('E38').Value = '  -0.08%  '
$ws.Range('E39').Value = '  +19.93%  '
$ws.Range('E40').Value = '  -1.89%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.42'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.99%  '
$ws.Range('B42').Value = 'MXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.19'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.784'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.28%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '63.38'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.61%  '
$ws.Range('D45').Value = '1.737.03'
$ws.Range('E45').Value = '  +0.21%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '91.16'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.59%  '
$ws.Range('E47').Value = '  -1.76%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₆0105'
$ws.Range('E48').Value = '  -1.00%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.101'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.02%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0510'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.64%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.00'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.27%  '
